$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 228-229, pushing the existing rows 228..239 down to 230..241.
$ws.Rows("228:229").Insert()

# New row 228 data
$ws.Cells.Item(228,1).Value  = 7
$ws.Cells.Item(228,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(228,3).Value  = "Ñuble"
$ws.Cells.Item(228,4).Value  = 44753
$ws.Cells.Item(228,5).Value  = 16
$ws.Cells.Item(228,6).Value  = 100112009
$ws.Cells.Item(228,7).Value  = "Acelga"
$ws.Cells.Item(228,8).Value  = "Sin especificar"
$ws.Cells.Item(228,9).Value  = "Primera"
$ws.Cells.Item(228,10).Value = 200
$ws.Cells.Item(228,11).Value = 600
$ws.Cells.Item(228,12).Value = 700
$ws.Cells.Item(228,13).Value = 650
$ws.Cells.Item(228,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(228,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(228,16).Value = 650
$ws.Cells.Item(228,17).Value = 1
$ws.Cells.Item(228,18).Value = "Hortaliza"

# New row 229 data
$ws.Cells.Item(229,1).Value  = 7
$ws.Cells.Item(229,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(229,3).Value  = "Ñuble"
$ws.Cells.Item(229,4).Value  = 44753
$ws.Cells.Item(229,5).Value  = 16
$ws.Cells.Item(229,6).Value  = 100112009
$ws.Cells.Item(229,7).Value  = "Acelga"
$ws.Cells.Item(229,8).Value  = "Sin especificar"
$ws.Cells.Item(229,9).Value  = "Segunda"
$ws.Cells.Item(229,10).Value = 100
$ws.Cells.Item(229,11).Value = 500
$ws.Cells.Item(229,12).Value = 500
$ws.Cells.Item(229,13).Value = 500
$ws.Cells.Item(229,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(229,15).Value = "Provincia de Diguillín"
$ws.Cells.Item(229,16).Value = 500
$ws.Cells.Item(229,17).Value = 1
$ws.Cells.Item(229,18).Value = "Hortaliza"
